$d = $word.ActiveDocument

# --- Change 1 (paragraph 4): "nummers" -> "liedjes" (2x) -----------------
$p4 = $d.Paragraphs(4)
if ($p4.Range.Text -like "*de nummers te gaan bijhouden. De nummers zullen*") {
    $p4.Range.Text = "We zorgen ervoor dat onze album klasse voorzien wordt van de correct properties om de titel, de artiest en de liedjes te gaan bijhouden. De liedjes zullen we voornu als String gaan opslaan in een array. Deze zal voor elk album anders zijn dus zal je op de correcte plek moeten initialiseren"
}

# --- Change 2 (paragraph 6): "nummers" -> "liedjes" -----------------------
$p6 = $d.Paragraphs(6)
if ($p6.Range.Text -like "*met de nummers van enkele*") {
    $p6.Range.Text = "Maak enkele stringArray's aan met de liedjes van enkele van je favoriete albums. "
}

# --- Change 3 (paragraph 9): "absolut" -> "absolute" ----------------------
$p9 = $d.Paragraphs(9)
if ($p9.Range.Text -like "*je absolut favorite singles*") {
    $p9.Range.Text = "Maak ook eem mix tape van je absolute favorite singles en vernoem de artiest naar jezelf.  "
}

# --- Change 4 (paragraph 11, first run only): "nummer" -> "liedje" --------
# This paragraph has multiple runs (trailing " ", "kijk dan naar Arrays.copyOff",
# "4 - MuziekFeestje"). Assigning only the replacement for the first sentence
# lets the engine diff it against the untouched remainder, leaving the other
# runs (and their formatting / the <w:br/>) completely intact.
$p11 = $d.Paragraphs(11)
if ($p11.Range.Text -like "*kunnen nummer toevoegen aan onze array*") {
    $p11.Range.Text = "We gaan nu ook onze Album klasse gaan uitbreiden zodat we kunnen liedje toevoegen aan onze array. Zoek hiervoor op hoe je best een array vergroot en een object toevoegd. Probeer eerst met array manipulatie zelf"
}
